$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in existing meeting date/time string (row 26 / B26)
$ws.Range("B26").Value = "10/31 /1:00"
$ws.Range("B26").NumberFormat = "General"

# Fill in the new attendance row (row 27) for the 11/3 meeting
$ws.Range("B27").Value = "11/3 /4:15"
$ws.Range("B27").NumberFormat = "General"
$ws.Range("C27").Value = "Google Hangout"
$ws.Range("D27").Value = "A"
$ws.Range("E27").Value = "U"
$ws.Range("F27").Value = "A"
$ws.Range("G27").Value = "U"
$ws.Range("H27").Value = "A"
$ws.Range("I27").Value = "A"

$ws.Range("C27").Font.Bold = $true
$ws.Range("C27").Borders.Item(8).LineStyle = -4142

# Meeting place filled in ahead of time for the following two rows
$ws.Range("C28").Value = "Google Hangout"
$ws.Range("C28").Font.Bold = $true
$ws.Range("C28").Borders.Item(8).LineStyle = -4142

$ws.Range("C29").Value = "Google Hangout"
$ws.Range("C29").Font.Bold = $true
$ws.Range("C29").Borders.Item(8).LineStyle = -4142

$ws.Range("I27").Select()
